# Updated cryptos list on Fri Apr 28 05:23:25 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') cells hold numeric-looking text (e.g. '1.000', '326.98')
# that must stay as literal text, not be coerced into real numbers. A leading
# apostrophe (Excel's standard text-entry prefix) forces that, same as typing
# it by hand; the apostrophe itself is not stored as part of the value.

# Row 2
$ws.Cells.Item(2, 4).Value = '''29.507.97'
$ws.Cells.Item(2, 5).Value = '  +1.44%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '''1.914.60'
$ws.Cells.Item(3, 5).Value = '  +0.24%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.000'
$ws.Cells.Item(4, 5).Value = '  +0.06%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''326.98'
$ws.Cells.Item(5, 5).Value = '  -1.95%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''1.0000'
$ws.Cells.Item(6, 5).Value = '  +0.03%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.4787'
$ws.Cells.Item(7, 5).Value = '  +3.10%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.4098'
$ws.Cells.Item(8, 5).Value = '  +0.03%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''47.78'
$ws.Cells.Item(9, 5).Value = '  +0.20%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.08030'
$ws.Cells.Item(10, 5).Value = '  +0.20%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.27%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''22.37'
$ws.Cells.Item(12, 5).Value = '  +2.38%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''1.876.70'
$ws.Cells.Item(13, 5).Value = '  -1.28%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''5.948'
$ws.Cells.Item(14, 5).Value = '  -0.05%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''7.146'
$ws.Cells.Item(15, 5).Value = '  +0.67%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''89.56'
$ws.Cells.Item(16, 5).Value = '  +0.42%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.04%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.06616'
$ws.Cells.Item(18, 5).Value = '  +0.70%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''0.00001032'
$ws.Cells.Item(19, 5).Value = '  +0.07%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''17.74'
$ws.Cells.Item(20, 5).Value = '  +1.13%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.09%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''29.522.55'
$ws.Cells.Item(22, 5).Value = '  +1.50%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''5.552'
$ws.Cells.Item(23, 5).Value = '  +1.91%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''11.51'
$ws.Cells.Item(24, 5).Value = '  +1.79%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -1.60%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''2.106.15'
$ws.Cells.Item(26, 5).Value = '  -1.02%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''153.69'
$ws.Cells.Item(27, 5).Value = '  -2.28%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''19.81'
$ws.Cells.Item(28, 5).Value = '  +0.22%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''5.792'
$ws.Cells.Item(29, 5).Value = '  +6.75%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''2.134'
$ws.Cells.Item(30, 5).Value = '  +0.84%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''117.63'
$ws.Cells.Item(31, 5).Value = '  -0.80%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''1.060'
$ws.Cells.Item(32, 5).Value = '  +7.03%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''0.09566'
$ws.Cells.Item(33, 5).Value = '  +1.53%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''1.422'
$ws.Cells.Item(34, 5).Value = '  -0.63%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''3.574'
$ws.Cells.Item(35, 5).Value = '  -0.60%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''5.389'
$ws.Cells.Item(36, 5).Value = '  +1.26%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.06104'
$ws.Cells.Item(37, 5).Value = '  -0.07%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''0.02252'
$ws.Cells.Item(38, 5).Value = '  +0.24%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''8.339'
$ws.Cells.Item(39, 5).Value = '  -0.60%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''1.174'
$ws.Cells.Item(40, 5).Value = '  -0.64%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.5882'
$ws.Cells.Item(41, 5).Value = '  +1.02%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.1844'
$ws.Cells.Item(42, 5).Value = '  +0.92%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.59%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Cronos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(44, 4).Value = '''0.08023'
$ws.Cells.Item(44, 5).Value = '  +13.77%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).Value = '''2.463'
$ws.Cells.Item(45, 5).Value = '  +4.89%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(46, 4).Value = '''1.292'
$ws.Cells.Item(46, 5).Value = '  +2.22%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '''0.5543'
$ws.Cells.Item(47, 5).Value = '  +0.48%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Value = '''12.14'
$ws.Cells.Item(48, 5).Value = '  +0.23%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''1.929'
$ws.Cells.Item(49, 5).Value = '  +0.47%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''113.62'
$ws.Cells.Item(50, 5).Value = '  +1.78%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''44.74'
$ws.Cells.Item(51, 5).Value = '  -6.44%  '
